# Replace abbreviated state names in column Q (rows 2-17) of each sheet
# with the corresponding full state name.

$map = @{
    "EUM"    = "Estados Unidos Mexicanos"
    "CDMX"   = "Ciudad de México"
    "Dgo."   = "Durango"
    "Gto."   = "Guanajuato"
    "Gro."   = "Guerrero"
    "Hgo."   = "Hidalgo"
    "Jal."   = "Jalisco"
    "Mex."   = "Estado de México"
    "Mich."  = "Michoacán"
    "Mor."   = "Morelos"
    "Nay."   = "Nayarit"
    "Ags."   = "Aguascalientes"
    "NL"     = "Nuevo León"
    "Oax."   = "Oaxaca"
    "Pue."   = "Puebla"
    "Qro."   = "Querétaro"
    "Q. Roo" = "Quintana Roo"
    "SLP"    = "San Luis Potosí"
    "Sin."   = "Sinaloa"
    "Son."   = "Sonora"
    "Tab."   = "Tabasco"
    "Tamps." = "Tamaulipas"
    "BC"     = "Baja California"
    "Tlax."  = "Tlaxcala"
    "Ver."   = "Veracruz"
    "Yuc."   = "Yucatán"
    "Zac."   = "Zacatecas"
    "BCS"    = "Baja California Sur"
    "Camp."  = "Campeche"
    "Coah."  = "Coahuila"
    "Col."   = "Colima"
    "Chis."  = "Chiapas"
    "Chih."  = "Chihuahua"
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($row = 2; $row -le 17; $row++) {
        $cell = $ws.Cells.Item($row, 17)  # Column Q is the 17th column
        $current = $cell.Text
        if ($map.ContainsKey($current)) {
            $cell.Value = $map[$current]
        }
    }
}
